# Update row-2 numeric results on each year sheet with fresh server values.
$wb = $excel.ActiveWorkbook

# Ordered list (one entry per worksheet, in sheet order: 2025, 2030, 2035, 2040, 2045, 2050)
$updates = @(
    @{ "A2" = 0;                  "B2" = 65.23935292841448;  "E2" = 28962.47207921735;
       "I2" = 26428.08358355595;  "L2" = 26903.0948845862;   "M2" = 11698.763646295;
       "N2" = 7232.873454107385;  "O2" = 6979.511843720443 },
    @{ "A2" = 144.8958715432552;  "B2" = 3653.322435707081;  "E2" = 45443.70038969377;
       "I2" = 47349.97900353095;  "L2" = 26903.0948845862;   "M2" = 17830.16682873526;
       "N2" = 9238.989105154791;  "O2" = 8048.129841190603 },
    @{ "A2" = 5087.980127543257;  "B2" = 6633.645507455078;  "E2" = 51090.70498697401;
       "I2" = 62699.13762201397;  "L2" = 26903.0948845862;   "M2" = 23809.18592562385;
       "N2" = 11559.8841274993;   "O2" = 13823.3599135127 },
    @{ "A2" = 5087.980127543257;  "B2" = 6633.645507455078;  "E2" = 51090.70498697401;
       "I2" = 62699.13762201397;  "L2" = 26903.0948845862;   "M2" = 23809.18592562385;
       "N2" = 11559.8841274993;   "O2" = 13823.3599135127 },
    @{ "A2" = 5087.980127543257;  "B2" = 6633.645507455078;  "E2" = 51090.70498697401;
       "I2" = 62699.13762201397;  "L2" = 26903.0948845862;   "M2" = 23809.18592562385;
       "N2" = 11559.8841274993;   "O2" = 13823.3599135127 },
    @{ "A2" = 5087.980127543257;  "B2" = 6633.645507455078;  "E2" = 51090.70498697401;
       "I2" = 62699.13762201397;  "L2" = 26903.0948845862;   "M2" = 23809.18592562385;
       "N2" = 11559.8841274993;   "O2" = 13823.3599135127 }
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $cellValues = $updates[$i]
    foreach ($addr in $cellValues.Keys) {
        $ws.Range($addr).Value = $cellValues[$addr]
    }
}
